$d = $word.ActiveDocument
$t = $d.Tables(1)

# 1) First paragraph ("Weekly Schedule:") spacing: w:after 0 -> 240 (twips) = 12pt
$d.Paragraphs(1).SpaceAfter = 12

# 2) Table borders: set all sides/insides to "none" (val=none sz=0 space=0 color=auto)
$t.Borders.Enable = $false
$b = $t.Borders(-1)
$b.LineWidth = 0
$b.ColorIndex = 0
$b.LineStyle = 0

# 3) Column widths: 1303/8232 dxa -> 1328/8207 dxa (dxa = points * 20)
$t.Columns(1).Width = 66.4
$t.Columns(2).Width = 410.35

# 4) Row heights: add trHeight val=1008 (twips) = 50.4 pt, for all 5 rows
$t.Rows.Height = 50.4

# 5) Day-name cells: append a new, separate run containing ":" after the day text
function Add-ColonRun([int]$row) {
    $cell = $t.Cell($row, 1)
    $full = $cell.Range
    $textRng = $d.Range($full.Start, $full.End - 1)
    $day = $textRng.Text
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $day + '</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $textRng.InsertXML($xml)
}

Add-ColonRun 1
Add-ColonRun 2
Add-ColonRun 3
Add-ColonRun 4
Add-ColonRun 5

# 6) Friday's content cell: collapse the three runs into a single run
$fridayCell = $t.Cell(4, 2)
$fullRng = $fridayCell.Range
$textRng2 = $d.Range($fullRng.Start, $fullRng.End - 1)
$combinedText = "Respond to another student" + [char]0x2019 + "s post on the discussion board by midnight on Friday. "
$xml2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">' + $combinedText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$textRng2.InsertXML($xml2)

Write-Host "Edit complete"
